# Generate Report for Handback
# Updates the localization-status report after a handback cycle completed:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#  - Latest Handback DateTime values are refreshed for zh-cn and de-de
#  - The stale "handback file not latest" error is cleared now that it's in sync
#  - Status/Error Detail columns are widened/narrowed to fit their new content

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ----- Overview sheet -----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ----- zh-cn sheet -----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-08-22 14:56:25"
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333334

# ----- de-de sheet -----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-08-22 14:56:32"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333334

Write-Host "Report regenerated for handback."
